$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "59.691.48" or "539.07")
# that must stay plain text, exactly like the source workbook. Force the
# cell to Text format while assigning, then restore the default "Normal"
# style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "59.691.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  +0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.615.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  +1.12%  "
$ws.Range("E4").Value2 = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "539.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "142.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +2.15%  "
$ws.Range("E7").Value2 = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "6.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  +0.73%  "
$ws.Range("E10").Value2 = "  +1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.336"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "3.072.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "59.601.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "20.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "2.629.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +1.89%  "
$ws.Range("E17").Value2 = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "341.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  +0.11%  "
$ws.Range("E19").Value2 = "  +1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "10.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "6.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -1.43%  "
$ws.Range("E22").Value2 = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "67.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.410"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +1.13%  "
$ws.Range("E25").Value2 = "  -1.19%  "
$ws.Range("E26").Value2 = "  +0.21%  "
$ws.Range("E27").Value2 = "  +2.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.0₃0748"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +3.50%  "
$ws.Range("E29").Value2 = "  +0.01%  "
$ws.Range("E30").Value2 = "  +5.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "5.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "18.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +0.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "150.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +1.08%  "
$ws.Range("E34").Value2 = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +0.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.839"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +3.15%  "
$ws.Range("E37").Value2 = "  -0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.830"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +0.37%  "
$ws.Range("E39").Value2 = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "278.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +2.97%  "
$ws.Range("E41").Value2 = "  -0.02%  "
$ws.Range("E42").Value2 = "  +0.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "10.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -0.15%  "
$ws.Range("E44").Value2 = "  -0.23%  "
$ws.Range("E45").Value2 = "  +2.00%  "
$ws.Range("B46").Value2 = "VeChain"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.0223"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +0.55%  "
$ws.Range("B47").Value2 = "Maker"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.945.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -1.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "18.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "4.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "110.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -3.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "4.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +0.94%  "
